$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: "short term guardian appointment form"
#      -> "short-term guardian appointment form"
# (hyphenate "short term" -> "short-term" in the heading sentence)
# ------------------------------------------------------------------
$r1 = $d.Content
$old1 = "short term guardian appointment form"
$new1 = "short-term guardian appointment form"
$found1 = $r1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "edit1 found=$found1"

# ------------------------------------------------------------------
# Edit 2: append a new sentence right after
# "Keep your copy in a safe place." in the "Make copies" list item.
# ------------------------------------------------------------------
$r2 = $d.Content
$old2 = "Keep your copy in a safe place."
$new2 = "Keep your copy in a safe place. You may want to keep it with your child’s birth certificate."
$found2 = $r2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Output "edit2 found=$found2"
